$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Fill in the new log entry on row 15 (last 4 digits, date, start time, end time, description)
$ws.Range("B15").Value = """0624"""
$ws.Range("C15").Value = 43925
$ws.Range("D15").Value = "7:25pm"
$ws.Range("E15").Value = "7:41pm"
$ws.Range("G15").Value = "Proofread the document once and made some changes, will submit tomorrow"

# Scroll the view over to column G and select G15, matching the author's final view state
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 7
$ws.Range("G15").Select()
